# Apply the "Coursera P4E" commit:
#  - Core-Courses sheet: row 17 (Study Designs in Epidemiology) gets Sr number 13
#  - Health-Care sheet: row 12 gets a new course ("Data Science in Stratified
#    Healthcare and Precision Medicine") with a hyperlink to its Coursera page
#  - Health-Care sheet: active selection moves from D12 to E12

$wb = $excel.ActiveWorkbook

# --- Core-Courses sheet -----------------------------------------------
$wsCore = $wb.Worksheets.Item("Core-Courses")
$wsCore.Range("B17").Value = 13

# --- Health-Care sheet --------------------------------------------------
$wsHealth = $wb.Worksheets.Item("Health-Care")

$wsHealth.Range("C12").Value = "Data Science in Stratified Healthcare and Precision Medicine"

$wsHealth.Hyperlinks.Add(
    $wsHealth.Range("D12"),
    "https://www.coursera.org/learn/datascimed",
    "",
    "",
    "https://www.coursera.org/learn/datascimed"
)

# Hyperlinks.Add re-applies its own ad-hoc "Hyperlink" formatting, which
# would otherwise leave D12 on a freshly minted style distinct from (but
# visually identical to) the one already shared by D5:D11. Re-copy the
# formatting from the row above so D12 keeps using that same existing
# hyperlink style, while leaving the newly created hyperlink relationship
# (and the display text it set) untouched.
$wsHealth.Range("D11").Copy()
$wsHealth.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to match the saved workbook state (E12).
$wsHealth.Activate()
$wsHealth.Range("E12").Select()
